# Daily attendance processing - 2025-12-21 13:55:49
#
# The "Recorded By" column (G) lists the users who recorded/updated each
# attendance row as a comma-separated string. Whenever the automated
# "System" account appears as the FIRST entry in that list, this script
# re-orders the list so "System" is no longer first (by reversing the
# full comma-separated sequence). Rows whose first entry is something
# other than "System" (e.g. a real user) are left untouched, as are rows
# that only contain a single value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($null -eq $value) {
        continue
    }

    $text = [string]$value
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -ge 2 -and $parts[0].ToLower() -eq "system") {
        $reversed = $parts[($parts.Length - 1)..0]
        $cell.Value2 = [string]::Join(", ", $reversed)
    }
}
